# TC09_Canine_Filter_StageOfDisease-5a.xlsx
#
# The "FilesTab" Neo4j query stored in B4 of the "startup" sheet is
# rewritten to drop the `File Type` and `Breed` output columns
# (two `coalesce(...)` lines removed from the RETURN clause).
# The previous selection (E2) is moved onto the edited cell (B4),
# matching the author re-saving the file right after editing that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFilesTabQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['Va']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Cells.Item(4, 2).Value = $newFilesTabQuery

$ws.Range("B4").Select() | Out-Null
